# names.xlsx - "some fixes with encoding"
# Rewrites the three person names in column A and gives the list a
# proper bordered / centered / word-wrapped Times New Roman look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Replace the three names (shared strings) with the new values.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Жуков Александр Аркадьевич"
$ws.Range("A2").Value = "Жуков Дмитрий Алексеевич"
$ws.Range("A3").Value = "Жуков Никита Юрьевич"

# ---------------------------------------------------------------------
# 2. Build the new cell format (font / border / alignment) once on a
#    scratch cell, then copy just the formatting onto A1:A3 in a single
#    paste-special. Doing it this way (instead of setting each property
#    straight on A1:A3) keeps the generated style table minimal - only
#    one new font/border/cell style is actually used by the sheet.
# ---------------------------------------------------------------------
$scratch = $ws.Range("Z100")
$scratch.Font.Name = "Times New Roman"
$scratch.Font.Color = 0
$scratch.Borders.LineStyle = 1
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4108
$scratch.WrapText = $true

$scratch.Copy()
$ws.Range("A1:A3").PasteSpecial(-4122)
$scratch.Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Move the active selection, matching the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("D13").Select()
